$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 930
$ws1.Range("F8").Value = 81
$ws1.Range("F9").Value = 575
$ws1.Range("F15").Value = 1069
$ws1.Range("F17").Value = 6772
$ws1.Range("F21").Value = 7697
$ws1.Range("F24").Value = 3866
$ws1.Range("F26").Value = 2201
$ws1.Range("F33").Value = 247
$ws1.Range("F36").Value = 1842
$ws1.Range("F38").Value = 213
$ws1.Range("F40").Value = 518
$ws1.Range("F42").Value = 1298
$ws1.Range("F44").Value = 1964
$ws1.Range("F45").Value = 2162
$ws1.Range("F46").Value = 11

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 10
$ws2.Range("F8").Value = 86

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 81
$ws4.Range("F9").Value = 575
$ws4.Range("F15").Value = 1069
$ws4.Range("F17").Value = 6772
$ws4.Range("F21").Value = 7697
$ws4.Range("F24").Value = 3866
$ws4.Range("F26").Value = 2201
$ws4.Range("F34").Value = 10
$ws4.Range("F35").Value = 247
$ws4.Range("F36").Value = 1842
$ws4.Range("F38").Value = 213
$ws4.Range("F40").Value = 518
$ws4.Range("F43").Value = 1298
$ws4.Range("F45").Value = 1964
$ws4.Range("F47").Value = 2162
$ws4.Range("F48").Value = 11
$ws4.Range("F49").Value = 86
